$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New response rows appended to the survey sheet (rows 582-600).
# Each entry lists the A-L column values plus either an M or N answer
# (the sheet uses mutually-exclusive M/N columns for a branching question).
$newRows = @(
    @{ A=45189.94790576389; B='hshs0104746@naver.com'; C='광고홍보학과'; D=20232639.0; E='최희수'; F='74:26'; G=0.2; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='952만 명'; J=0.059; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Red'; LastCol="M"; LastVal='모름/무응답' }
    @{ A=45189.97729207176; B='ehdus040127@naver.com'; C='사회복지학부'; D=20232307.0; E='김도연'; F='74:26'; G=0.2; H='OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다.'; I='166만 명'; J=0.002; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Black'; LastCol="N"; LastVal='모름/무응답' }
    @{ A=45190.02308297454; B='hjkiubb@naver.com'; C='미디어스쿨'; D=20232538.0; E='박재은'; F='74:26'; G=0.2; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='952만 명'; J=0.059; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Black'; LastCol="N"; LastVal='국민부담률을 OECD 평균 수준으로 높여야 한다' }
    @{ A=45190.03219640046; B='smile001118@naver.com'; C='사회학과'; D=20222240.0; E='홍성준'; F='74:26'; G=0.2; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='952만 명'; J=0.059; K='상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다'; L='Red'; LastCol="M"; LastVal='국민부담률을 OECD 평균 수준으로 높여야 한다' }
    @{ A=45190.072204861106; B='poliku8630@naver.com'; C='컨텐츠 IT'; D=20205197.0; E='심지혁'; F='76:24'; G=0.2; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='952만 명'; J=0.059; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Red'; LastCol="M"; LastVal='국민부담률을 OECD 평균 수준으로 높여야 한다' }
    @{ A=45190.07256223379; B='hyeeun7356@gmail.com'; C='식품영양학과'; D=20203824.0; E='유혜은'; F='74:26'; G=0.2; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='952만 명'; J=0.059; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Black'; LastCol="N"; LastVal='모름/무응답' }
    @{ A=45190.14519309028; B='dnjsgmlwjd1020@naver.com'; C='인문학부'; D=20231057.0; E='원희정'; F='74:26'; G=0.2; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='952만 명'; J=0.059; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Black'; LastCol="N"; LastVal='국민부담률을 OECD 평균 수준으로 높여야 한다' }
    @{ A=45190.384077395836; B='minjeong7432@gmail.com'; C='간호학과'; D=20236217.0; E='김민정'; F='74:26'; G=0.2; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='952만 명'; J=0.059; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Black'; LastCol="N"; LastVal='국민부담률을 아일랜드 수준으로 낮춰야 한다' }
    @{ A=45190.41397086806; B='alsgk03@naver.com'; C='사회학과'; D=20222213.0; E='박민하'; F='74:26'; G=0.2; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='952만 명'; J=0.059; K='중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'; L='Black'; LastCol="N"; LastVal='모름/무응답' }
    @{ A=45190.52351451389; B='suani3176@gmail.com'; C='사회복지학부'; D=20232328.0; E='박수안'; F='75:25'; G=0.2; H='프랑스와 스웨덴의 국민부담률은 꾸준히 40%를 넘고 있다.'; I='779만 명'; J=0.374; K='중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'; L='Black'; LastCol="N"; LastVal='국민부담률을 OECD 평균 수준으로 높여야 한다' }
    @{ A=45190.588201342594; B='hkmcosmos1@gmail.com'; C='글로벌 비즈니스'; D=20226429.0; E='한기민'; F='77:23'; G=0.15; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='779만 명'; J=0.151; K='중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'; L='Black'; LastCol="N"; LastVal='국민부담률을 아일랜드 수준으로 낮춰야 한다' }
    @{ A=45190.61006407408; B='yeon-jin22@naver.com'; C='데이터사이언스학부'; D=20233256.0; E='최연진'; F='74:26'; G=0.2; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='952만 명'; J=0.059; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Black'; LastCol="N"; LastVal='국민부담률을 아일랜드 수준으로 낮춰야 한다' }
    @{ A=45190.62784372685; B='bsw030409@naver.com'; C='철학과'; D=20221043.0; E='백승우'; F='74:26'; G=0.2; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='952만 명'; J=0.059; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Black'; LastCol="N"; LastVal='모름/무응답' }
    @{ A=45190.63231763889; B='hyj4213@naver.com'; C='미디어스쿨'; D=20232590.0; E='함영준'; F='74:26'; G=0.2; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='952만 명'; J=0.059; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Red'; LastCol="M"; LastVal='국민부담률을 아일랜드 수준으로 낮춰야 한다' }
    @{ A=45190.69257645833; B='seo1020102p@naver.com'; C='미디어스쿨'; D=20232537.0; E='박재연'; F='74:26'; G=0.2; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='952만 명'; J=0.059; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Red'; LastCol="M"; LastVal='국민부담률을 아일랜드 수준으로 낮춰야 한다' }
    @{ A=45190.70098991898; B='kddong99@gmail.com'; C='빅데이터전공'; D=20181205.0; E='김상준'; F='76:24'; G=0.2; H='OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다.'; I='952만 명'; J=0.151; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Red'; LastCol="M"; LastVal='국민부담률을 OECD 평균 수준으로 높여야 한다' }
    @{ A=45190.750975717594; B='psjj3840@gmail.com'; C='디지털미디어콘텐츠'; D=20215154.0; E='박서진'; F='74:26'; G=0.2; H='우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다.'; I='779만 명'; J=0.151; K='상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다'; L='Red'; LastCol="M"; LastVal='모름/무응답' }
    @{ A=45190.78272258102; B='ub030801@naver.com'; C='간호학과'; D=20226256.0; E='신유빈'; F='77:23'; G=0.15; H='OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다.'; I='779만 명'; J=0.374; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Red'; LastCol="M"; LastVal='모름/무응답' }
    @{ A=45190.78992811343; B='yenaridia@naver.com'; C='경영학과'; D=20201634.0; E='최예나'; F='74:26'; G=0.2; H='미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'; I='952만 명'; J=0.059; K='법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'; L='Red'; LastCol="M"; LastVal='국민부담률을 아일랜드 수준으로 낮춰야 한다' }
)

$startRow = 582
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("A" + $r).Value = $row.A
    $ws.Range("B" + $r).Value = $row.B
    $ws.Range("C" + $r).Value = $row.C
    $ws.Range("D" + $r).Value = $row.D
    $ws.Range("E" + $r).Value = $row.E
    $ws.Range("F" + $r).Value = $row.F
    $ws.Range("G" + $r).Value = $row.G
    $ws.Range("H" + $r).Value = $row.H
    $ws.Range("I" + $r).Value = $row.I
    $ws.Range("J" + $r).Value = $row.J
    $ws.Range("K" + $r).Value = $row.K
    $ws.Range("L" + $r).Value = $row.L
    $ws.Range($row.LastCol + $r).Value = $row.LastVal
}

# Copy the cell formatting (number formats / styles) from the last pre-existing
# data rows so the new rows end up with the same style indices as the source:
#   row 581 supplies styles for columns A-L and for the "M" answer column,
#   row 580 supplies the style for the "N" answer column (581 has no N cell).
$srcBase = $ws.Range("A581:L581")
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $dstBase = $ws.Range("A" + $r + ":L" + $r)
    $srcBase.Copy() | Out-Null
    $dstBase.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $row = $newRows[$i]
    if ($row.LastCol -eq "M") {
        $srcLast = $ws.Range("M581")
    } else {
        $srcLast = $ws.Range("N580")
    }
    $dstLast = $ws.Range($row.LastCol + $r)
    $srcLast.Copy() | Out-Null
    $dstLast.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0
Write-Host "Added rows 582-600 to sheet"